$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep zoom level consistent (matches source change of zoomScale/zoomScaleNormal)
$excel.ActiveWindow.Zoom = 225

# Fill in new row 5 data: Name, First name, Date
$ws.Range("A5").Value = "anao "
$ws.Range("B5").Value = "gerard "
$ws.Range("C5").Value = (Get-Date -Year 2023 -Month 6 -Day 3 -Hour 0 -Minute 0 -Second 0)

# Copy the date number format from C4 to C5 (numFmtId 14)
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update selection to A5
$ws.Range("A5").Select()
